$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 74 with revised figures ---
$ws.Range("B74").Value = 2556
$ws.Range("C74").Value = 439
$ws.Range("D74").Value = -355
$ws.Range("E74").Value = 795
$ws.Range("F74").Value = 64
$ws.Range("H74").Value = 62
$ws.Range("I74").Value = 3696
$ws.Range("J74").Value = 1484
$ws.Range("L74").Value = 1134
$ws.Range("Q74").Value = -783
$ws.Range("R74").Value = 910
$ws.Range("S74").Value = 1646
$ws.Range("T74").Value = 1778
$ws.Range("V74").Value = 1811

# --- Add new row 75 with the 01-04-2021 data ---
# Enter the period label via a text formula so Excel does not
# auto-convert the dd-mm-yyyy-looking text into a date value, then
# convert the formula result back into a plain (shared-string) value.
$ws.Range("A75").Formula = '="01-04-2021"'
$ws.Range("A75").Copy()
$ws.Range("A75").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("B75").Value = 4476
$ws.Range("C75").Value = 10932
$ws.Range("D75").Value = 7769
$ws.Range("E75").Value = 3163
$ws.Range("F75").Value = 28
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 28
$ws.Range("I75").Value = 1142
$ws.Range("J75").Value = 891
$ws.Range("K75").Value = 562
$ws.Range("L75").Value = -311
$ws.Range("M75").Value = -8347
$ws.Range("N75").Value = -533
$ws.Range("O75").Value = -520
$ws.Range("P75").Value = -13
$ws.Range("Q75").Value = 1255
$ws.Range("R75").Value = 1913
$ws.Range("S75").Value = 2562
$ws.Range("T75").Value = 2517
$ws.Range("U75").Value = -148
$ws.Range("V75").Value = 2665
$ws.Range("W75").Value = 45
